# Regional Availability Factor workbook update
# - Updates the "last updated" date on the About sheet
# - Sets RAF-capacity demand-altering-tech rows (hydrogen turbines) to full capacity credit (1)
# - Leaves the workbook with the RAF-capacity sheet active/selected, matching the saved view state

$wb = $excel.ActiveWorkbook

# --- About sheet: bump the revision date (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- RAF-generation sheet: no longer the tab shown when the file is opened ---
$wsGen = $wb.Worksheets.Item("RAF-generation")
$wsGen.Range("A1").RowHeight = 29.5

# --- RAF-demand-altering-techs sheet: matching row-1 height touch-up ---
$wsDemand = $wb.Worksheets.Item("RAF-demand-altering-techs")
$wsDemand.Range("A1").RowHeight = 29.5

# --- RAF-capacity sheet: update capacity credit multipliers for hydrogen plants ---
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1
$wsCap.Range("A:A").ColumnWidth = 28.209635416666668

# Make RAF-capacity the active/selected sheet and restore its saved view
$wsCap.Activate() | Out-Null
$wsCap.Range("B25").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 80
$win.ScrollRow = 14
$win.ScrollColumn = 1
